$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 45.82032266666667
$ws.Range("H2").Value = 137.460968
$ws.Range("I2").Value = 0.2798800263398926
$ws.Range("J2").Value = 0.2798800263398927
$ws.Range("M2").Value = 16.28832166666666
$ws.Range("N2").Value = 48.864965
$ws.Range("O2").Value = 0.2220374022124247
$ws.Range("P2").Value = 0.2220374022124247
$ws.Range("Q2").Value = 746.3361544651244
$ws.Range("R2").Value = 6717.02539018612
$ws.Range("S2").Value = 0.06214383397965475
$ws.Range("T2").Value = 0.06214383397965476
$ws.Range("G3").Value = 45.82032266666667
$ws.Range("H3").Value = 137.460968
$ws.Range("I3").Value = 0.2798800263398926
$ws.Range("J3").Value = 0.2798800263398927
$ws.Range("N3").Value = 65.980473
$ws.Range("O3").Value = 0.299808519696413
$ws.Range("P3").Value = 0.2998085196964129
$ws.Range("Q3").Value = 1007.74885418643
$ws.Range("R3").Value = 9069.739687677866
$ws.Range("S3").Value = 0.08391041638955628
$ws.Range("T3").Value = 0.08391041638955628
$ws.Range("G4").Value = 45.82032266666667
$ws.Range("H4").Value = 137.460968
$ws.Range("I4").Value = 0.2798800263398926
$ws.Range("J4").Value = 0.2798800263398927
$ws.Range("M4").Value = 19.34010233333333
$ws.Range("N4").Value = 58.020307
$ws.Range("O4").Value = 0.2636383396948583
$ws.Range("P4").Value = 0.2636383396948583
$ws.Range("Q4").Value = 886.1697293196863
$ws.Range("R4").Value = 7975.527563877177
$ws.Range("S4").Value = 0.07378710545800249
$ws.Range("T4").Value = 0.0737871054580025
$ws.Range("G5").Value = 45.82032266666667
$ws.Range("H5").Value = 137.460968
$ws.Range("I5").Value = 0.2798800263398926
$ws.Range("J5").Value = 0.2798800263398927
$ws.Range("M5").Value = 15.736544
$ws.Range("N5").Value = 47.209632
$ws.Range("O5").Value = 0.2145157383963041
$ws.Range("P5").Value = 0.2145157383963041
$ws.Range("Q5").Value = 721.0535237381974
$ws.Range("R5").Value = 6489.481713643776
$ws.Range("S5").Value = 0.06003867051267911
$ws.Range("T5").Value = 0.06003867051267911
$ws.Range("I6").Value = 0.2427284602664133
$ws.Range("J6").Value = 0.2427284602664133
$ws.Range("M6").Value = 16.28832166666666
$ws.Range("N6").Value = 48.864965
$ws.Range("O6").Value = 0.2220374022124247
$ws.Range("P6").Value = 0.2220374022124247
$ws.Range("Q6").Value = 647.2667163267822
$ws.Range("R6").Value = 5825.40044694104
$ws.Range("S6").Value = 0.05389479676057615
$ws.Range("T6").Value = 0.05389479676057615
$ws.Range("I7").Value = 0.2427284602664133
$ws.Range("J7").Value = 0.2427284602664133
$ws.Range("N7").Value = 65.980473
$ws.Range("O7").Value = 0.299808519696413
$ws.Range("P7").Value = 0.2998085196964129
$ws.Range("Q7").Value = 873.9792221358989
$ws.Range("R7").Value = 7865.812999223089
$ws.Range("S7").Value = 0.07277206036066296
$ws.Range("T7").Value = 0.07277206036066294
$ws.Range("I8").Value = 0.2427284602664133
$ws.Range("J8").Value = 0.2427284602664133
$ws.Range("M8").Value = 19.34010233333333
$ws.Range("N8").Value = 58.020307
$ws.Range("O8").Value = 0.2636383396948583
$ws.Range("P8").Value = 0.2636383396948583
$ws.Range("Q8").Value = 768.5386368773992
$ws.Range("R8").Value = 6916.847731896592
$ws.Range("S8").Value = 0.06399252826132656
$ws.Range("T8").Value = 0.06399252826132658
$ws.Range("I9").Value = 0.2427284602664133
$ws.Range("J9").Value = 0.2427284602664133
$ws.Range("M9").Value = 15.736544
$ws.Range("N9").Value = 47.209632
$ws.Range("O9").Value = 0.2145157383963041
$ws.Range("P9").Value = 0.2145157383963041
$ws.Range("Q9").Value = 625.3401283237547
$ws.Range("R9").Value = 5628.061154913792
$ws.Range("S9").Value = 0.05206907488384761
$ws.Range("T9").Value = 0.05206907488384761
$ws.Range("G10").Value = 43.41682666666667
$ws.Range("H10").Value = 130.25048
$ws.Range("I10").Value = 0.2651989746877358
$ws.Range("J10").Value = 0.2651989746877358
$ws.Range("M10").Value = 16.28832166666666
$ws.Range("N10").Value = 48.864965
$ws.Range("O10").Value = 0.2220374022124247
$ws.Range("P10").Value = 0.2220374022124247
$ws.Range("Q10").Value = 707.1872384925778
$ws.Range("R10").Value = 6364.6851464332
$ws.Range("S10").Value = 0.05888409140906342
$ws.Range("T10").Value = 0.05888409140906344
$ws.Range("G11").Value = 43.41682666666667
$ws.Range("H11").Value = 130.25048
$ws.Range("I11").Value = 0.2651989746877358
$ws.Range("J11").Value = 0.2651989746877358
$ws.Range("N11").Value = 65.980473
$ws.Range("O11").Value = 0.299808519696413
$ws.Range("P11").Value = 0.2998085196964129
$ws.Range("Q11").Value = 954.8875865418936
$ws.Range("R11").Value = 8593.988278877041
$ws.Range("S11").Value = 0.07950891202613657
$ws.Range("T11").Value = 0.07950891202613657
$ws.Range("G12").Value = 43.41682666666667
$ws.Range("H12").Value = 130.25048
$ws.Range("I12").Value = 0.2651989746877358
$ws.Range("J12").Value = 0.2651989746877358
$ws.Range("M12").Value = 19.34010233333333
$ws.Range("N12").Value = 58.020307
$ws.Range("O12").Value = 0.2636383396948583
$ws.Range("P12").Value = 0.2636383396948583
$ws.Range("Q12").Value = 839.6858707219291
$ws.Range("R12").Value = 7557.172836497361
$ws.Range("S12").Value = 0.0699166173754534
$ws.Range("T12").Value = 0.06991661737545342
$ws.Range("G13").Value = 43.41682666666667
$ws.Range("H13").Value = 130.25048
$ws.Range("I13").Value = 0.2651989746877358
$ws.Range("J13").Value = 0.2651989746877358
$ws.Range("M13").Value = 15.736544
$ws.Range("N13").Value = 47.209632
$ws.Range("O13").Value = 0.2145157383963041
$ws.Range("P13").Value = 0.2145157383963041
$ws.Range("Q13").Value = 683.2308031803734
$ws.Range("R13").Value = 6149.07722862336
$ws.Range("S13").Value = 0.0568893538770824
$ws.Range("T13").Value = 0.05688935387708241
$ws.Range("G14").Value = 34.73892266666667
$ws.Range("H14").Value = 104.216768
$ws.Range("I14").Value = 0.2121925387059582
$ws.Range("J14").Value = 0.2121925387059582
$ws.Range("M14").Value = 16.28832166666666
$ws.Range("N14").Value = 48.864965
$ws.Range("O14").Value = 0.2220374022124247
$ws.Range("P14").Value = 0.2220374022124247
$ws.Range("Q14").Value = 565.8387467481244
$ws.Range("R14").Value = 5092.54872073312
$ws.Range("S14").Value = 0.04711468006313033
$ws.Range("T14").Value = 0.04711468006313034
$ws.Range("G15").Value = 34.73892266666667
$ws.Range("H15").Value = 104.216768
$ws.Range("I15").Value = 0.2121925387059582
$ws.Range("J15").Value = 0.2121925387059582
$ws.Range("N15").Value = 65.980473
$ws.Range("O15").Value = 0.299808519696413
$ws.Range("P15").Value = 0.2998085196964129
$ws.Range("Q15").Value = 764.0301830190294
$ws.Range("R15").Value = 6876.271647171265
$ws.Range("S15").Value = 0.06361713092005714
$ws.Range("T15").Value = 0.06361713092005714
$ws.Range("G16").Value = 34.73892266666667
$ws.Range("H16").Value = 104.216768
$ws.Range("I16").Value = 0.2121925387059582
$ws.Range("J16").Value = 0.2121925387059582
$ws.Range("M16").Value = 19.34010233333333
$ws.Range("N16").Value = 58.020307
$ws.Range("O16").Value = 0.2636383396948583
$ws.Range("P16").Value = 0.2636383396948583
$ws.Range("Q16").Value = 671.8543193230863
$ws.Range("R16").Value = 6046.688873907777
$ws.Range("S16").Value = 0.05594208860007577
$ws.Range("T16").Value = 0.05594208860007577
$ws.Range("G17").Value = 34.73892266666667
$ws.Range("H17").Value = 104.216768
$ws.Range("I17").Value = 0.2121925387059582
$ws.Range("J17").Value = 0.2121925387059582
$ws.Range("M17").Value = 15.736544
$ws.Range("N17").Value = 47.209632
$ws.Range("O17").Value = 0.2145157383963041
$ws.Range("P17").Value = 0.2145157383963041
$ws.Range("Q17").Value = 546.6705850565974
$ws.Range("R17").Value = 4920.035265509376
$ws.Range("S17").Value = 0.04551863912269496
$ws.Range("T17").Value = 0.04551863912269496
